# Generate Report for Handoff
# The file "b.md" has been re-handed-off for both zh-cn and de-de locales.
# Update status from "Handed back: in sync with en-US" to "Ready for handoff"
# on the Overview sheet and on each locale sheet, and record the new
# handoff file name + datetime for the b.md row on each locale sheet.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# --- Overview sheet: row for b.md (row 3) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# --- zh-cn sheet: row for b.md (row 3) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = $newStatus
$zhcn.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("D3").Value = "2016-02-19 08:00:03"

foreach ($h in $zhcn.Hyperlinks) {
    $rngAddr = $h.Range.Address()
    if ($rngAddr -eq '$C$3') {
        $h.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
    }
}

# --- de-de sheet: row for b.md (row 3) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = $newStatus
$dede.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("D3").Value = "2016-02-19 08:00:20"

foreach ($h in $dede.Hyperlinks) {
    $rngAddr = $h.Range.Address()
    if ($rngAddr -eq '$C$3') {
        $h.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
    }
}
